# participation_ratios.xlsx update: add IBM pad-width / corner-radius ratio
# sweeps to the Voltages sheet, and fix the Capacitances C84 formula.

$wb = $excel.ActiveWorkbook
$wsV = $wb.Worksheets.Item("Voltages")
$wsC = $wb.Worksheets.Item("Capacitances")

# --- Voltages: pad-width sweep header (row 58) ---
$wsV.Range("A58").Value = "Pad Width"
$wsV.Range("A58").Font.Bold = $true
$wsV.Range("A58").NumberFormat = "0.00E+00"

# --- Voltages: pad-width sweep rows 59-64 ---
$wsV.Range("A59").Value = "IBM_ratio_pad_w_300"
$wsV.Range("B59").Value = 14.1416
$wsV.Range("C59").Value = 6.5128
$wsV.Range("D59").Value = [double]"1.4e-08"
$wsV.Range("E59").Value = [double]"4.27261e-14"
$wsV.Range("F59").Value = [double]"4.712e-05"
$wsV.Range("G59").Value = 0.00035077
$wsV.Range("H59").Value = 0.00010817
$wsV.Range("J59").Value = 300
$wsV.Range("K59").Value = 18
$wsV.Range("L59").Value = 60
$wsV.Range("M59").Value = 650
$wsV.Range("N59").Value = "-"
$wsV.Range("O59").Value = 735711

$wsV.Range("A60").Value = "IBM_ratio_pad_w_300_openboundary"
$wsV.Range("B60").Value = 14.1415
$wsV.Range("C60").Value = 6.5172
$wsV.Range("D60").Value = [double]"1.4e-08"
$wsV.Range("E60").Value = [double]"4.27261e-14"
$wsV.Range("J60").Value = 300
$wsV.Range("K60").Value = 18
$wsV.Range("L60").Value = 60
$wsV.Range("M60").Value = 650
$wsV.Range("N60").Value = "-"
$wsV.Range("O60").Value = 830555

$wsV.Range("A61").Value = "IBM_ratio_pad_w_350"
$wsV.Range("B61").Value = 14.1414
$wsV.Range("C61").Value = 6.51
$wsV.Range("D61").Value = [double]"1.22e-08"
$wsV.Range("E61").Value = [double]"4.908e-14"
$wsV.Range("F61").Value = [double]"4.2855e-05"
$wsV.Range("G61").Value = 0.00034443
$wsV.Range("H61").Value = 0.00010417
$wsV.Range("J61").Value = 350
$wsV.Range("K61").Value = 18
$wsV.Range("L61").Value = 60
$wsV.Range("M61").Value = 650
$wsV.Range("N61").Value = "-"
$wsV.Range("O61").Value = 832072

$wsV.Range("A62").Value = "IBM_ratio_pad_w_400"
$wsV.Range("B62").Value = 14.1413
$wsV.Range("C62").Value = 6.5096
$wsV.Range("D62").Value = [double]"1.08e-08"
$wsV.Range("E62").Value = [double]"5.54507e-14"
$wsV.Range("F62").Value = [double]"4.2694e-05"
$wsV.Range("G62").Value = 0.00033775
$wsV.Range("H62").Value = 0.00010182
$wsV.Range("J62").Value = 400
$wsV.Range("K62").Value = 18
$wsV.Range("L62").Value = 60
$wsV.Range("M62").Value = 650
$wsV.Range("N62").Value = "-"
$wsV.Range("O62").Value = 935539

$wsV.Range("A63").Value = "IBM_ratio_pad_w_450"
$wsV.Range("B63").Value = 14.1412
$wsV.Range("C63").Value = 6.5072
$wsV.Range("D63").Value = [double]"9.69e-09"
$wsV.Range("E63").Value = [double]"6.18554e-14"
$wsV.Range("F63").Value = [double]"4.1739e-05"
$wsV.Range("G63").Value = 0.00033197
$wsV.Range("H63").Value = 0.00010068
$wsV.Range("J63").Value = 450
$wsV.Range("K63").Value = 18
$wsV.Range("L63").Value = 60
$wsV.Range("M63").Value = 650
$wsV.Range("N63").Value = "-"
$wsV.Range("O63").Value = 1038493

$wsV.Range("A64").Value = "IBM_ratio_pad_w_500"
$wsV.Range("B64").Value = 14.1411
$wsV.Range("C64").Value = 6.5052
$wsV.Range("D64").Value = [double]"8.78e-09"
$wsV.Range("E64").Value = [double]"6.83157e-14"
$wsV.Range("F64").Value = [double]"4.0369e-05"
$wsV.Range("G64").Value = 0.0003281
$wsV.Range("H64").Value = [double]"9.9211e-05"
$wsV.Range("J64").Value = 500
$wsV.Range("K64").Value = 18
$wsV.Range("L64").Value = 60
$wsV.Range("M64").Value = 650
$wsV.Range("N64").Value = "-"
$wsV.Range("O64").Value = 1146237

$wsV.Range("I59").Formula = "=SUM(F59:H59)"
$wsV.Range("I60:I64").Formula = "=SUM(F60:H60)"

$wsV.Range("A59:A64").NumberFormat = "0.00E+00"
$wsV.Range("D59:D65").NumberFormat = "0.00E+00"
$wsV.Range("E59:E64").NumberFormat = "0.00E+00"
$wsV.Range("F59:F64").NumberFormat = "0.00E+00"
$wsV.Range("G59:G64").NumberFormat = "0.00E+00"
$wsV.Range("H59:H64").NumberFormat = "0.00E+00"
$wsV.Range("I59:I64").NumberFormat = "0.00E+00"
$wsV.Range("J59:J64").Style = "40% - Accent1"
$wsV.Range("A65").NumberFormat = "0.00E+00"

# --- Voltages: corner-radius sweep header (row 66) ---
$wsV.Range("A66").Value = "Corner Radius"
$wsV.Range("A66").Font.Bold = $true

# --- Voltages: corner-radius sweep rows 67-72 ---
$wsV.Range("A67").Value = "IBM_ratio_w450_l60_cornerr_5"
$wsV.Range("B67").Value = 14.1413
$wsV.Range("C67").Value = 6.5048
$wsV.Range("D67").Value = [double]"9.71e-09"
$wsV.Range("E67").Value = [double]"6.172e-14"
$wsV.Range("F67").Value = [double]"3.7383e-05"
$wsV.Range("G67").Value = 0.00030276
$wsV.Range("H67").Value = [double]"9.5348e-05"
$wsV.Range("J67").Value = 450
$wsV.Range("K67").Value = 18
$wsV.Range("L67").Value = 60
$wsV.Range("M67").Value = 650
$wsV.Range("N67").Value = 5
$wsV.Range("O67").Value = 835079

$wsV.Range("A68").Value = "IBM_ratio_w450_l60_cornerr_10"
$wsV.Range("B68").Value = 14.1413
$wsV.Range("C68").Value = 6.5044
$wsV.Range("D68").Value = [double]"9.75e-09"
$wsV.Range("E68").Value = [double]"6.1464e-14"
$wsV.Range("F68").Value = [double]"3.557e-05"
$wsV.Range("G68").Value = 0.00030078
$wsV.Range("H68").Value = [double]"9.4668e-05"
$wsV.Range("J68").Value = 450
$wsV.Range("K68").Value = 18
$wsV.Range("L68").Value = 60
$wsV.Range("M68").Value = 650
$wsV.Range("N68").Value = 10
$wsV.Range("O68").Value = 831573

$wsV.Range("A69").Value = "IBM_ratio_w450_l60_cornerr_15"
$wsV.Range("B69").Value = 14.1413
$wsV.Range("C69").Value = 6.5028
$wsV.Range("D69").Value = [double]"9.81e-09"
$wsV.Range("E69").Value = [double]"6.11259e-14"
$wsV.Range("F69").Value = [double]"3.6656e-05"
$wsV.Range("G69").Value = 0.00029846
$wsV.Range("H69").Value = [double]"9.6023e-05"
$wsV.Range("J69").Value = 450
$wsV.Range("K69").Value = 18
$wsV.Range("L69").Value = 60
$wsV.Range("M69").Value = 650
$wsV.Range("N69").Value = 15
$wsV.Range("O69").Value = 824442

$wsV.Range("A70").Value = "IBM_ratio_w450_l60_cornerr_20"
$wsV.Range("B70").Value = 14.1414
$wsV.Range("C70").Value = 6.502
$wsV.Range("D70").Value = [double]"9.88e-09"
$wsV.Range("E70").Value = [double]"6.0706e-14"
$wsV.Range("F70").Value = [double]"3.3695e-05"
$wsV.Range("G70").Value = 0.00029976
$wsV.Range("H70").Value = [double]"9.5086e-05"
$wsV.Range("J70").Value = 450
$wsV.Range("K70").Value = 18
$wsV.Range("L70").Value = 60
$wsV.Range("M70").Value = 650
$wsV.Range("N70").Value = 20
$wsV.Range("O70").Value = 817046

$wsV.Range("A71").Value = "IBM_ratio_w450_l60_cornerr_25"
$wsV.Range("B71").Value = 14.1413
$wsV.Range("C71").Value = 6.5256
$wsV.Range("D71").Value = [double]"9.96e-09"
$wsV.Range("E71").Value = [double]"6.0221e-14"
$wsV.Range("F71").Value = [double]"3.4653e-05"
$wsV.Range("G71").Value = 0.00029683
$wsV.Range("H71").Value = [double]"9.4099e-05"
$wsV.Range("J71").Value = 450
$wsV.Range("K71").Value = 18
$wsV.Range("L71").Value = 60
$wsV.Range("M71").Value = 650
$wsV.Range("N71").Value = 25
$wsV.Range("O71").Value = 799697

$wsV.Range("A72").Value = "IBM_ratio_w450_l60_cornerr_29"
$wsV.Range("B72").Value = 14.1413
$wsV.Range("C72").Value = 6.5124
$wsV.Range("D72").Value = [double]"1e-08"
$wsV.Range("E72").Value = [double]"5.97734e-14"
$wsV.Range("F72").Value = [double]"3.4888e-05"
$wsV.Range("G72").Value = 0.00029882
$wsV.Range("H72").Value = [double]"9.4731e-05"
$wsV.Range("J72").Value = 450
$wsV.Range("K72").Value = 18
$wsV.Range("L72").Value = 60
$wsV.Range("M72").Value = 650
$wsV.Range("N72").Value = 29
$wsV.Range("O72").Value = 799697

$wsV.Range("I67").Formula = "=SUM(F67:H67)"
$wsV.Range("I68:I72").Formula = "=SUM(F68:H68)"

$wsV.Range("A67:A72").NumberFormat = "0.00E+00"
$wsV.Range("D67:D72").NumberFormat = "0.00E+00"
$wsV.Range("E67:E72").NumberFormat = "0.00E+00"
$wsV.Range("F67:F72").NumberFormat = "0.00E+00"
$wsV.Range("G67:G72").NumberFormat = "0.00E+00"
$wsV.Range("H67:H72").NumberFormat = "0.00E+00"
$wsV.Range("I67:I72").NumberFormat = "0.00E+00"
$wsV.Range("N67:N72").Style = "40% - Accent1"

# --- Capacitances: C84 breaks out of the C67:C92 shared formula ---
$wsC.Range("C84").Formula = '=1/((2*PI()*$B$1)^2*B84)'

# --- view state + page setup ---
$wsV.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$wsV.Range("D37").Select()
$wsV.Range("A16").Select()

$wsC.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsC.Range("B20").Select()

$ps = $wsC.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

$wsV.Activate()
